# Adding colors to first row in excel output
# Also inserts a "Department" column (new column E) with value
# "T301 - Green House Controlled" for the Greenhouse budget rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert a new column before column E (shifts E.. rightward by one).
#    Column D ("Role Ending" header / "418 - Plant Growth" data) itself
#    is left untouched by the insert.
# ---------------------------------------------------------------------
$ws.Columns("E:E").Insert()

# New column header for the inserted column.
$ws.Range("E2").Value = "Department"

# Rows 3-14: move the legacy "418 - Plant Growth" text that is still
# sitting in D into the new E column, and give D its new department
# value.
for ($r = 3; $r -le 14; $r++) {
    $ws.Cells.Item($r, 5).Value = "418 - Plant Growth"
    $ws.Cells.Item($r, 4).Value = "T301 - Green House Controlled"
}

# ---------------------------------------------------------------------
# 2. Color the project-header cells on row 1 (now shifted one column
#    to the right because of the insert above).
# ---------------------------------------------------------------------
$ws.Range("L1:Q1").Interior.Color = 13223074    # 00A2C4C9 - B10/20 AgPlenus
$ws.Range("R1:V1").Interior.Color = 13421812    # 00F4CCCC - B20/20 Lavie-Bio (1)
$ws.Range("W1:Z1").Interior.Color = 13882323    # 00D3D3D3 - B20/20 Lavie-Bio (2)
$ws.Range("AA1:AD1").Interior.Color = 10275833  # 00F9CB9C - B40/20 CPB
$ws.Range("AE1:AJ1").Interior.Color = 13888217  # 00D9EAD3 - B70/20 Biomica + B74/20 Canonic (1)
$ws.Range("AK1:AO1").Interior.Color = 13431551  # 00FFF2CC - B74/20 Canonic (2)
$ws.Range("AP1:AY1").Interior.Color = 15254943  # 009FC5E8 - B80/20 PRoduct
$ws.Range("AZ1:BE1").Interior.Color = 8242323   # 0093C47D - B72/20 Casterra
$ws.Range("BF1:BG1").Interior.Color = 14471658  # 00EAD1DC - trailing empty cells
